# Oct 28th - Status
# Adds the day's status rows (28/10/2021) to the daily status tracker sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A61").Value = "28/10/2021"
$ws.Range("B61").Value = "We divided the yavta testapp as I analyzed the code and the teammate tried on the output debug"

$ws.Range("B62").Value = "Studied each struct in yavta testapp listed 16 structs"
$ws.Range("C62").Value = "log files of yavta pushed to the GitHub ,need to clarify the doubts"

$ws.Range("B63").Value = "Explored bytes perline"
$ws.Range("C63").Value = "Exploring the v4l2 headers in yavta testapp"

[void]$ws.Range("B64").Select()
